# Update the workbook to reflect the new "MIT battery cost decline" source
# (replacing the old BNEF "New Energy Outlook 2018" source) and switch the
# headline number on the PDiBCpDoC sheet from a hard-coded value to a
# formula that averages the two learning-rate bounds quoted in the source.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("PDiBCpDoC")

# --- "About" sheet -------------------------------------------------------

# Drop the old "last updated" date stamp in C1 (no longer tracked).
$ws1.Range("C1").Clear()

# New source citation.
$ws1.Range("B3").Value = "Massachusetts Institute of Technology"
$ws1.Range("B4").Value = 2021
$ws1.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$ws1.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$ws1.Range("B7").Value = "Abstract"

# Old note (graph only extends to 2030) no longer applies - clear it but
# keep the italic-style cell in place.
$ws1.Range("C8").Value = $null

# New footnote explaining how the figure was derived.
$ws1.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# Remove the embedded chart picture that illustrated the old BNEF source.
[void]$ws1.Shapes.Item(1).Delete()

# --- "PDiBCpDoC" sheet ----------------------------------------------------

# Replace the hard-coded 18% with the average of the quoted 20%-27% range.
$ws2.Range("B2").Formula = "=AVERAGE(0.2,0.27)"

[void]$ws2.Range("I4").Select()

# Leave "About" as the active sheet/tab, with its selection moved to A10.
[void]$ws1.Range("A10").Select()
